# Reduce the magnitude of disturbances.
# The "f" column (C) holds frequency values oscillating around 1.0.
# Shrink the deviation from 1.0 by a factor of 10 for rows 3-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shrink the deviation of the "f" (frequency) column from 1.0 by a
# factor of 10, i.e. new = 1 + (old - 1) / 10, for the disturbance rows.
$newValues = @{
    3  = 1.001
    4  = 1.002
    5  = 1.003
    6  = 1.002
    7  = 1.001
    9  = 0.999
    10 = 0.998
    11 = 0.997
    12 = 0.998
    13 = 0.999
}

foreach ($r in $newValues.Keys) {
    $ws.Cells.Item($r, 3).Value2 = $newValues[$r]
}

# Update the view state to match: scrolled down one row, selection on C14.
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("C14").Select()
